$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the end time for row 35 (E35) from 12:00 to 11:45
$ws.Range("E35").Value = 0.48958333333333331

# Move the active selection to E36
$ws.Range("E36").Select()
